$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.983.04"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.742.03"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "250.93"
$ws.Range("E5").Value = "  +6.95%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.5159"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "1.742.79"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "0.07231"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "15.14"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "0.6510"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "4.636"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "77.72"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "25.999.70"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "0.000006828"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "1.966.27"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "4.279"
$ws.Range("D23").Value = "8.686"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "5.374"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").Value = "136.23"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D27").Value = "15.26"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "1.782"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "105.94"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "3.969"
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("D31").Value = "0.08228"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "3.655"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "0.04708"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").Value = "2.657"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "0.9972"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "0.6241"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "1.920"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "1.000"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "100.54"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "0.3860"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "6.320"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "0.1133"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").Value = "55.68"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").Value = "0.05229"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "30.79"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "7.530"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "0.3433"
$ws.Range("E51").Value = "  -1.08%  "
